$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '59.094.08'
$ws.Range('E2').Value = '  -1.64%  '
$ws.Range('D3').Value = '2.330.17'
$ws.Range('E3').Value = '  -3.58%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '551.63'
$ws.Range('E5').Value = '  -0.32%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '132.21'
$ws.Range('E6').Value = '  -3.30%  '
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.574'
$ws.Range('E8').Value = '  -2.30%  '
$ws.Range('D9').Value = '2.330.33'
$ws.Range('E9').Value = '  -3.47%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.103'
$ws.Range('E10').Value = '  -2.73%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '5.58'
$ws.Range('E11').Value = '  -2.31%  '
$ws.Range('E12').Value = '  +1.11%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '0.337'
$ws.Range('E13').Value = '  -4.89%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '24.21'
$ws.Range('E14').Value = '  -2.28%  '
$ws.Range('D15').Value = '2.740.75'
$ws.Range('E15').Value = '  -3.79%  '
$ws.Range('D16').Value = '59.049.48'
$ws.Range('E16').Value = '  -1.58%  '
$ws.Range('E17').Value = '  -2.47%  '
$ws.Range('D18').Value = '2.247.97'
$ws.Range('E18').Value = '  -7.09%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '10.81'
$ws.Range('E19').Value = '  -3.91%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '4.36'
$ws.Range('E20').Value = '  -3.07%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '317.05'
$ws.Range('E21').Value = '  -3.18%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '6.54'
$ws.Range('E22').Value = '  -3.78%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '0.999'
$ws.Range('E23').Value = '  -0.06%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '63.67'
$ws.Range('E24').Value = '  -1.98%  '
$ws.Range('E25').Value = '  -2.71%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '0.999'
$ws.Range('E26').Value = '  -0.56%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '8.12'
$ws.Range('E27').Value = '  -5.68%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '1.32'
$ws.Range('E28').Value = '  -6.11%  '
$ws.Range('E29').Value = '  +0.51%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '169.89'
$ws.Range('E30').Value = '  -0.28%  '
$ws.Range('D31').Value = '0.0₃0735'
$ws.Range('E31').Value = '  -4.80%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '1.13'
$ws.Range('E32').Value = '  +5.21%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '5.87'
$ws.Range('E33').Value = '  -3.79%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.386'
$ws.Range('E34').Value = '  -3.78%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.998'
$ws.Range('E35').Value = '  -0.07%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '17.87'
$ws.Range('E36').Value = '  -3.46%  '
$ws.Range('E37').Value = '  +0.07%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '1.28'
$ws.Range('E38').Value = '  -4.33%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '4.02'
$ws.Range('E39').Value = '  -4.64%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '38.21'
$ws.Range('E40').Value = '  -1.45%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '1.53'
$ws.Range('E41').Value = '  -3.98%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '305.96'
$ws.Range('E42').Value = '  -5.21%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '143.08'
$ws.Range('E43').Value = '  -1.71%  '
$ws.Range('E44').Value = '  -4.99%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.0956'
$ws.Range('E45').Value = '  -0.66%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.0505'
$ws.Range('E46').Value = '  -1.95%  '
$ws.Range('B47').Value = 'Mantle'
$ws.Range('C47').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.565'
$ws.Range('E47').Value = '  -2.18%  '
$ws.Range('B48').Value = 'InjectiveProtocol'
$ws.Range('C48').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '18.83'
$ws.Range('E48').Value = '  -4.84%  '
$ws.Range('E49').Value = '  -2.27%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '16.86'
$ws.Range('E50').Value = '  -3.37%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '11.03'
$ws.Range('E51').Value = '  -0.19%  '
